# Fill in the daily-update rows for 9/19/2021 (serial date 44458) that were
# previously left blank, and move the active-cell selection down to A22 —
# matching the next blank row after the newly filled-in entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: Self-learning entry, 7 hours
$ws.Range("A20").Value = 44458
$ws.Range("A20").NumberFormat = "d-mmm"
$ws.Range("B20").Value = "7 hours"
$ws.Range("C20").Value = "Self learning: Watched a couple of crash courses on asp.net mvc and continued with the pluralsight course"

# Row 21: Task entry, 1 hour
$ws.Range("A21").Value = 44458
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("B21").Value = "1 hour"
$ws.Range("C21").Value = "Task: worked on the end user documentataion"

# Move selection to the next empty row, as recorded in the saved view state.
$ws.Range("A22").Select()
